# FreeCRMData.xlsx edit: complete the "ContactSheet2" data provider sheet
# with the full set of columns (matching "ContactSheet") and add a second
# data row, per the commit "Completed with Data provider with two rows."

$wb = $excel.ActiveWorkbook

$wsSrc  = $wb.Worksheets.Item("ContactSheet")
$wsDest = $wb.Worksheets.Item("ContactSheet2")

# ---------------------------------------------------------------------
# 1) ContactSheet2: extend header row (row 1) from column E to X, copying
#    the header labels already used on the "ContactSheet" sheet.
#    NOTE: ".Value" must be invoked as a method (parens) to read back the
#    real cell contents in this host - a bare property read returns the
#    .NET member-info string instead of the value.
# ---------------------------------------------------------------------
for ($col = 5; $col -le 24; $col++) {
    $srcVal = $wsSrc.Cells.Item(1, $col).Value()
    $wsDest.Cells.Item(1, $col).Value = $srcVal
}

# ---------------------------------------------------------------------
# 2) ContactSheet2: extend the first data row (row 2) from column E to X,
#    copying the sample values used on "ContactSheet" row 2.
# ---------------------------------------------------------------------
for ($col = 5; $col -le 24; $col++) {
    $srcVal = $wsSrc.Cells.Item(2, $col).Value()
    $wsDest.Cells.Item(2, $col).Value = $srcVal
}

# ---------------------------------------------------------------------
# 3) ContactSheet2: add a second data row (row 3) - same as row 2, except
#    FirstName/LastName/MiddleName get a distinct "...1" suffix.
# ---------------------------------------------------------------------
$wsDest.Cells.Item(3, 1).Value = "Latha1 "
$wsDest.Cells.Item(3, 2).Value = "Singh1"
$wsDest.Cells.Item(3, 3).Value = "Sri1"
for ($col = 4; $col -le 24; $col++) {
    $rowVal = $wsDest.Cells.Item(2, $col).Value()
    $wsDest.Cells.Item(3, $col).Value = $rowVal
}

# ---------------------------------------------------------------------
# 4) Hyperlinks for row 2 (Channel Link / Email Address / Personal Email)
#    - order matters: it drives the rId numbering (J,E,F).
# ---------------------------------------------------------------------
$wsDest.Hyperlinks.Add($wsDest.Range("J2"), "https://in.linkedin.com/")
$wsDest.Hyperlinks.Add($wsDest.Range("E2"), "mailto:abc@gmail.com")
$wsDest.Hyperlinks.Add($wsDest.Range("F2"), "mailto:abc@gmail.com")

# ... and for row 3.
$wsDest.Hyperlinks.Add($wsDest.Range("J3"), "https://in.linkedin.com/")
$wsDest.Hyperlinks.Add($wsDest.Range("E3"), "mailto:abc@gmail.com")
$wsDest.Hyperlinks.Add($wsDest.Range("F3"), "mailto:abc@gmail.com")

# Adding a Hyperlink re-styles the cell with a fresh style index; put the
# shared "Hyperlink" cell style back on all six cells so they reuse the
# same style already used elsewhere in the workbook.
$wsDest.Range("E2").Style = "Hyperlink"
$wsDest.Range("F2").Style = "Hyperlink"
$wsDest.Range("J2").Style = "Hyperlink"
$wsDest.Range("E3").Style = "Hyperlink"
$wsDest.Range("F3").Style = "Hyperlink"
$wsDest.Range("J3").Style = "Hyperlink"

# ---------------------------------------------------------------------
# 5) Selections: ContactSheet now has its whole data block selected,
#    ContactSheet2 selection moves to C3 (the last-typed cell).
# ---------------------------------------------------------------------
$wsSrc.Activate()
$wsSrc.Range("R1:X2").Select()

$wsDest.Activate()
$wsDest.Range("C3").Select()
